$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:49:58"
$ws1.Range("A3").Value = "Total filas: 4"

# Row 6
$ws1.Range("A6").Value = "02:49:58"
$ws1.Range("B6").Value = "03:01"
$ws1.Range("D6").Value = 12

# Row 7
$ws1.Range("A7").Value = "02:49:58"
$ws1.Range("D7").Value = 59

# Row 8
$ws1.Range("A8").Value = "02:49:58"
$ws1.Range("B8").Value = "04:01"
$ws1.Range("D8").Value = 72

# New row 9
$ws1.Range("A9").Value = "02:49:58"
$ws1.Range("B9").Value = "04:47"
$ws1.Range("C9").Value = "81_EL PELIGRO"
$ws1.Range("D9").Value = 118
$ws1.Range("E9").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 02:49:58"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 02:49:58"
